# Fixed a bug on debut and end date in add month page
# Adds a new bug-report row (row 4) describing that the categories / pie
# chart don't display on the "Single-month/details" page.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 4)
$ws.Cells.Item(4, 1).Value = 44493
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy"

$ws.Cells.Item(4, 2).Value = 2

$ws.Cells.Item(4, 3).Value = "Single-month/details"

$ws.Cells.Item(4, 4).Value = "Les catégories et le camembert ne s'affiche pas"

$ws.Cells.Item(4, 5).Value = "OUI"
$ws.Cells.Item(4, 5).HorizontalAlignment = -4108
$ws.Cells.Item(4, 5).VerticalAlignment = -4108

$ws.Cells.Item(4, 6).Value = "NON"
$ws.Cells.Item(4, 6).HorizontalAlignment = -4108
$ws.Cells.Item(4, 6).VerticalAlignment = -4108

# Header row (E2/F2) picks up the same centered style used by data cells
$ws.Cells.Item(2, 5).HorizontalAlignment = -4108
$ws.Cells.Item(2, 5).VerticalAlignment = -4108
$ws.Cells.Item(2, 6).HorizontalAlignment = -4108
$ws.Cells.Item(2, 6).VerticalAlignment = -4108

# Move the active selection as recorded after the edit
$ws.Range("D10").Select()
